$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the variable data columns (D, I, J, K, L, M, P) of rows 123..248 down
# by one row (into 124..249), working from the bottom up so we never
# overwrite a value before it has been copied.
$cols = @(4, 9, 10, 11, 12, 13, 16)
for ($r = 248; $r -ge 123; $r--) {
    foreach ($c in $cols) {
        $ws.Cells.Item($r + 1, $c).Value2 = $ws.Cells.Item($r, $c).Value2
    }
}

# Row 249 is brand new - it needs the rest of the (constant-across-the-block)
# columns copied over too, plus the date-number formatting on column D.
$last = 249
$prev = $last - 1
$ws.Cells.Item($last, 1).Value2  = $ws.Cells.Item($prev, 1).Value2
$ws.Cells.Item($last, 2).Value2  = $ws.Cells.Item($prev, 2).Value2
$ws.Cells.Item($last, 3).Value2  = $ws.Cells.Item($prev, 3).Value2
$ws.Cells.Item($last, 5).Value2  = $ws.Cells.Item($prev, 5).Value2
$ws.Cells.Item($last, 6).Value2  = $ws.Cells.Item($prev, 6).Value2
$ws.Cells.Item($last, 7).Value2  = $ws.Cells.Item($prev, 7).Value2
$ws.Cells.Item($last, 8).Value2  = $ws.Cells.Item($prev, 8).Value2
$ws.Cells.Item($last, 14).Value2 = $ws.Cells.Item($prev, 14).Value2
$ws.Cells.Item($last, 15).Value2 = $ws.Cells.Item($prev, 15).Value2
$ws.Cells.Item($last, 17).Value2 = $ws.Cells.Item($prev, 17).Value2
$ws.Cells.Item($last, 18).Value2 = $ws.Cells.Item($prev, 18).Value2
$ws.Cells.Item($last, 4).NumberFormat = $ws.Cells.Item($prev, 4).NumberFormat

# Finally, write the new weekly report row into row 123 (everything else in
# that row - A,B,C,E,F,G,H,I,N,O,Q,R - is unchanged from before the shift).
$ws.Cells.Item(123, 4).Value2  = 44966
$ws.Cells.Item(123, 10).Value2 = 35
$ws.Cells.Item(123, 11).Value2 = 40000
$ws.Cells.Item(123, 12).Value2 = 40000
$ws.Cells.Item(123, 13).Value2 = 40000
$ws.Cells.Item(123, 16).Value2 = 3077
